$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Planilha2" - the main baseACM table sheet

# --- Column A (type) was an unused placeholder column; the blank cells in
#     A2:A6 are removed entirely (clearing both value and formatting so the
#     <c> elements disappear from the XML). ---
$ws.Range("A2:A6").Clear() | Out-Null

# --- Normalize formatting on the author/title columns that previously
#     carried a redundant "blank" style (now unified with the default
#     style used elsewhere on the sheet). ---
$ws.Range("B2:B6").WrapText = $false
$ws.Range("C2").WrapText = $false
$ws.Range("C5").WrapText = $false
$ws.Range("G2").WrapText = $false

# --- Re-apply the wrapped-text style to the long title cells that keep it. ---
$ws.Range("C3").WrapText = $true
$ws.Range("C4").WrapText = $true
$ws.Range("C6").WrapText = $true

# --- Row 2 (Line Pouchard / "Revisiting the Data Lifecycle..."): mark the
#     paper as accepted (Status - Etapa 1) -- it already lists ACM (G2) as
#     the publisher. ---
$ws.Range("F2").Value = "Aceito"

# --- Rows 3-5 already show "Eliminado" in F; just drop the now-unused
#     blank publisher placeholder cells in column G. ---
$ws.Range("G3").Clear() | Out-Null
$ws.Range("G4").Clear() | Out-Null
$ws.Range("G5").Clear() | Out-Null

# --- Row 6 (Cerys Willoughby / "Encouraging and Facilitating..."): mark as
#     accepted too, and drop its blank publisher placeholder. ---
$ws.Range("F6").Value = "Aceito"
$ws.Range("G6").Clear() | Out-Null

# --- Update the saved selection / scroll position of the sheet. ---
$ws.Range("F7").Select() | Out-Null

Write-Host "done"
